$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.758.10'
$ws.Range("E2").Value = '  +0.77%  '

# Row 3
$ws.Range("D3").Value = '3.043.31'
$ws.Range("E3").Value = '  +0.75%  '

# Row 4
$ws.Range("E4").Value = '  -0.24%  '

# Row 5
$ws.Range("D5").Value = '''544.51'
$ws.Range("E5").Value = '  +1.67%  '

# Row 6
$ws.Range("D6").Value = '''133.97'
$ws.Range("E6").Value = '  +1.46%  '

# Row 7
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.28%  '

# Row 8
$ws.Range("D8").Value = '3.039.73'
$ws.Range("E8").Value = '  +1.03%  '

# Row 9
$ws.Range("D9").Value = '''0.495'
$ws.Range("E9").Value = '  +1.17%  '

# Row 10
$ws.Range("D10").Value = '''6.18'
$ws.Range("E10").Value = '  +1.43%  '

# Row 11
$ws.Range("D11").Value = '''0.148'
$ws.Range("E11").Value = '  -2.06%  '

# Row 12
$ws.Range("D12").Value = '''0.447'
$ws.Range("E12").Value = '  +0.83%  '

# Row 13
$ws.Range("D13").Value = '''0.0000223'
$ws.Range("E13").Value = '  +2.04%  '

# Row 14
$ws.Range("D14").Value = '''34.30'
$ws.Range("E14").Value = '  +2.21%  '

# Row 15
$ws.Range("D15").Value = '3.530.07'
$ws.Range("E15").Value = '  +0.07%  '

# Row 16
$ws.Range("D16").Value = '62.748.43'
$ws.Range("E16").Value = '  +0.47%  '

# Row 17
$ws.Range("D17").Value = '3.041.29'
$ws.Range("E17").Value = '  +0.01%  '

# Row 18
$ws.Range("E18").Value = '  -3.09%  '

# Row 19
$ws.Range("D19").Value = '''6.63'
$ws.Range("E19").Value = '  +1.69%  '

# Row 20
$ws.Range("D20").Value = '''477.62'
$ws.Range("E20").Value = '  +3.35%  '

# Row 21
$ws.Range("D21").Value = '''13.31'
$ws.Range("E21").Value = '  +0.82%  '

# Row 22
$ws.Range("D22").Value = '''0.673'
$ws.Range("E22").Value = '  -0.69%  '

# Row 23
$ws.Range("D23").Value = '''7.09'
$ws.Range("E23").Value = '  +3.04%  '

# Row 24
$ws.Range("D24").Value = '''80.98'
$ws.Range("E24").Value = '  +4.20%  '

# Row 25
$ws.Range("D25").Value = '''12.12'
$ws.Range("E25").Value = '  +1.86%  '

# Row 26
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  +0.33%  '

# Row 27
$ws.Range("D27").Value = '''2.71'
$ws.Range("E27").Value = '  +1.37%  '

# Row 28
$ws.Range("D28").Value = '''7.82'
$ws.Range("E28").Value = '  +2.24%  '

# Row 29
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.18%  '

# Row 30
$ws.Range("D30").Value = '''1.94'
$ws.Range("E30").Value = '  +5.70%  '

# Row 31
$ws.Range("D31").Value = '''25.72'
$ws.Range("E31").Value = '  +0.52%  '

# Row 32
$ws.Range("E32").Value = '  +0.02%  '

# Row 33
$ws.Range("D33").Value = '''2.38'
$ws.Range("E33").Value = '  +5.86%  '

# Row 34
$ws.Range("D34").Value = '''5.65'
$ws.Range("E34").Value = '  +6.03%  '

# Row 35
$ws.Range("D35").Value = '''54.96'
$ws.Range("E35").Value = '  -2.97%  '

# Row 36
$ws.Range("D36").Value = '''5.87'
$ws.Range("E36").Value = '  +0.66%  '

# Row 37
$ws.Range("D37").Value = '''462.33'
$ws.Range("E37").Value = '  +1.45%  '

# Row 38
$ws.Range("D38").Value = '3.159.70'
$ws.Range("E38").Value = '  -0.38%  '

# Row 39
$ws.Range("D39").Value = '''0.0800'
$ws.Range("E39").Value = '  +2.39%  '

# Row 40
$ws.Range("D40").Value = '''0.0390'
$ws.Range("E40").Value = '  +0.97%  '

# Row 41
$ws.Range("D41").Value = '''0.119'
$ws.Range("E41").Value = '  +2.36%  '

# Row 42
$ws.Range("D42").Value = '''8.10'
$ws.Range("E42").Value = '  +1.26%  '

# Row 43
$ws.Range("D43").Value = '''2.43'
$ws.Range("E43").Value = '  -1.04%  '

# Row 44
$ws.Range("D44").Value = '''26.79'
$ws.Range("E44").Value = '  +7.10%  '

# Row 46
$ws.Range("E46").Value = '  +0.12%  '

# Row 47
$ws.Range("E47").Value = '  +3.19%  '

# Row 48
$ws.Range("D48").Value = '''0.109'
$ws.Range("E48").Value = '  +0.97%  '

# Row 49
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '0.0₃0502'
$ws.Range("E49").Value = '  -0.73%  '

# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '''114.28'
$ws.Range("E50").Value = '  -5.35%  '

# Row 51
$ws.Range("E51").Value = '  +2.52%  '
